$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_Detail")

# --- Update the PO_Detail data row (Quote/DA?/UnitNoToMaintain/In Service Date) ---

# A2: Quote number 382245 -> 382425 (text, quote-prefixed "Text" number format)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "'382425"

# B2: DA? 2 -> 1 (plain text, default/"Normal" style)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

# C2: UnitNoToMaintain 00999159 -> 00999166 (keep leading zeros as text)
$ws.Range("C2").Value = "'00999166"

# D2: In Service Date cleared out (value removed, formatting kept)
$ws.Range("D2").Value = ""

# --- Move the active/selected sheet from Unit_to_Reconcile_Output to PO_Detail ---
$ws.Activate()
